# Insert a new quarterly column before column D, shifting existing data
# (old D..K) right into E..L, then populate the new column D with the
# newest quarter's figures (period ending 2018-09-30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column at D; this shifts D:K -> E:L and gives the new
#    D column the plain/default column style.
$ws.Range("D1").EntireColumn.Insert()

# 2) Copy the (now shifted) former-D column's formatting (still intact in
#    column E) into the new column D so the new cells carry the same
#    number formats (date format in row 7/38/80, thousands format
#    elsewhere) as the rest of the data.
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)

# 3) Populate the new column D with the latest quarter's values.
$ws.Range("D7").Value = 43373
$ws.Range("D8").Value = 33500
$ws.Range("D9").Value = 15900
$ws.Range("D10").Value = 17600
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = "NA"
$ws.Range("D15").Value = 6900
$ws.Range("D17").Value = 27700
$ws.Range("D18").Value = 5800
$ws.Range("D20").Value = 33700
$ws.Range("D21").Value = 46400
$ws.Range("D22").Value = 15600
$ws.Range("D23").Value = 23900
$ws.Range("D24").Value = 800
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 23100
$ws.Range("D27").Value = 22000
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -33700
$ws.Range("D33").Value = 22000
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 22000
$ws.Range("D38").Value = 43373
$ws.Range("D41").Value = 23800
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 127800
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 85000
$ws.Range("D48").Value = 1053500
$ws.Range("D49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 72700
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 1422900
$ws.Range("D57").Value = 35700
$ws.Range("D58").Value = 651000
$ws.Range("D59").Value = 0
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 475200
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 1205200
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -50600
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 217700
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43373
$ws.Range("D81").Value = 22000
$ws.Range("D83").Value = 6900
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = -21400
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -12800
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 37900
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 3700

# Row 91 ("Capital Expenditures") was re-keyed in full rather than being a
# pure shift of the old data, so set D91:L91 explicitly.
$ws.Range("D91").Value = -1000
$ws.Range("E91").Value = -2100
$ws.Range("F91").Value = -600
$ws.Range("G91").Value = -12000
$ws.Range("H91").Value = -300
$ws.Range("I91").Value = 400
$ws.Range("J91").Value = -1000
$ws.Range("K91").Value = 6300
$ws.Range("L91").Value = -6700
